{"js": "// Update master to output generated at c986bee\n// Replace the two-digit multiplication expressions in the table cells.\nconst replacements = [\n  [\"89\u00d766=\", \"60\u00d795=\"],\n  [\"71\u00d780=\", \"70\u00d722=\"],\n  [\"82\u00d733=\", \"42\u00d790=\"],\n  [\"23\u00d776=\", \"98\u00d789=\"],\n  [\"41\u00d745=\", \"30\u00d749=\"],\n  [\"70\u00d741=\", \"85\u00d782=\"],\n  [\"15\u00d711=\", \"49\u00d722=\"],\n  [\"98\u00d752=\", \"11\u00d782=\"],\n  [\"33\u00d765=\", \"44\u00d725=\"],\n  [\"73\u00d726=\", \"57\u00d779=\"],\n  [\"22\u00d712=\", \"35\u00d754=\"],\n  [\"61\u00d737=\", \"12\u00d727=\"],\n  [\"13\u00d766=\", \"93\u00d759=\"],\n  [\"76\u00d790=\", \"14\u00d754=\"],\n  [\"65\u00d768=\", \"56\u00d719=\"],\n  [\"97\u00d759=\", \"34\u00d730=\"],\n  [\"89\u00d773=\", \"56\u00d778=\"],\n  [\"70\u00d792=\", \"35\u00d724=\"],\n  [\"42\u00d772=\", \"13\u00d792=\"],\n  [\"16\u00d795=\", \"55\u00d726=\"],\n  [\"97\u00d750=\", \"24\u00d720=\"],\n  [\"87\u00d767=\", \"44\u00d792=\"],\n  [\"74\u00d752=\", \"39\u00d784=\"],\n  [\"50\u00d753=\", \"75\u00d739=\"],\n  [\"53\u00d788=\", \"82\u00d792=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update master to output generated at c986bee\n# Replace the two-digit multiplication expressions in the table cells.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"89\u00d766=\"; New = \"60\u00d795=\" },\n    @{ Old = \"71\u00d780=\"; New = \"70\u00d722=\" },\n    @{ Old = \"82\u00d733=\"; New = \"42\u00d790=\" },\n    @{ Old = \"23\u00d776=\"; New = \"98\u00d789=\" },\n    @{ Old = \"41\u00d745=\"; New = \"30\u00d749=\" },\n    @{ Old = \"70\u00d741=\"; New = \"85\u00d782=\" },\n    @{ Old = \"15\u00d711=\"; New = \"49\u00d722=\" },\n    @{ Old = \"98\u00d752=\"; New = \"11\u00d782=\" },\n    @{ Old = \"33\u00d765=\"; New = \"44\u00d725=\" },\n    @{ Old = \"73\u00d726=\"; New = \"57\u00d779=\" },\n    @{ Old = \"22\u00d712=\"; New = \"35\u00d754=\" },\n    @{ Old = \"61\u00d737=\"; New = \"12\u00d727=\" },\n    @{ Old = \"13\u00d766=\"; New = \"93\u00d759=\" },\n    @{ Old = \"76\u00d790=\"; New = \"14\u00d754=\" },\n    @{ Old = \"65\u00d768=\"; New = \"56\u00d719=\" },\n    @{ Old = \"97\u00d759=\"; New = \"34\u00d730=\" },\n    @{ Old = \"89\u00d773=\"; New = \"56\u00d778=\" },\n    @{ Old = \"70\u00d792=\"; New = \"35\u00d724=\" },\n    @{ Old = \"42\u00d772=\"; New = \"13\u00d792=\" },\n    @{ Old = \"16\u00d795=\"; New = \"55\u00d726=\" },\n    @{ Old = \"97\u00d750=\"; New = \"24\u00d720=\" },\n    @{ Old = \"87\u00d767=\"; New = \"44\u00d792=\" },\n    @{ Old = \"74\u00d752=\"; New = \"39\u00d784=\" },\n    @{ Old = \"50\u00d753=\"; New = \"75\u00d739=\" },\n    @{ Old = \"53\u00d788=\"; New = \"82\u00d792=\" }\n)\n\nforeach ($r in $replacements) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $r.Old\n    $find.Replacement.Text = $r.New\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($r.Old, $false, $false, $false, $false, $false, $true, 1, $false, $r.New, 2) | Out-Null\n}\n"}
